$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(1).RowHeight = 48
$ws.Rows.Item(2).RowHeight = 35.25
$ws.Rows.Item(8).RowHeight = 20.25
$ws.Rows.Item(17).RowHeight = 20.25
$ws.Rows.Item(19).RowHeight = 20.25
